$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A26").Value = "2023-12-06 14:01:43"
$ws.Range("B26").Value = 0.0006000000000000001
